# "modificata slide ex 1" - nudge the screenshot picture on the
# "CABIN MODULE" slide (slide 6) up a touch so it sits better under
# the code screenshot above it.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(6)
$sh = $s.Shapes.Item("Segnaposto contenuto 13")

# Move the picture placeholder up by ~23150 EMU (≈0.018"):
# from y=425769 EMU (33.52512 pt) to y=402619 EMU (31.70228 pt).
$sh.Top = 402619 / 914400 * 72
